# 8.9.1.1 - add a new "2020" data column (Q) mirroring the existing
# year header (row 4) and percentage value (row 5) formatting that is
# already used by column P.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Year header cell Q4 (matches D4:P4 style) ---
$ws.Range("Q4").Value = 2020
$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)   # xlPasteFormats

# --- Percentage value cell Q5 (matches D5:P5 style) ---
$ws.Range("Q5").Value = 3.3
$ws.Range("P5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)   # xlPasteFormats

# Clear the marching-ants clipboard marquee left over from the copies.
$excel.CutCopyMode = 0

# Match the author's recorded selection after the edit.
$null = $ws.Range("R4").Select()
